# Agregar "-1" para los registros "SinTurno" (Turno_Lista / ID_Turno)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "SinTurno"
$ws.Range("B26").Value = "-1"

# Reflect the view state from the authored edit: window scrolled so row 4
# is at the top, with the new row selected.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A26").Select() | Out-Null
